# Commit #2 // Pre-Alpha -> Clases y su respectivo Diagrama
#
# The last paragraph of the document ("Imágenes generadas por IA Copilot")
# gets its trailing word "Copilot" split out into its own (spell-check
# flagged) run, and a brand-new paragraph is appended after it with the
# "El motivo por el cual escojo este evento, ..." text.

$d = $word.ActiveDocument

# --- locate the paragraph that still reads the old, un-split sentence ---
$oldText = "Imágenes generadas por IA Copilot"
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq $oldText) {
        $target = $para
    }
}
if ($target -eq $null) {
    throw "Could not find the 'Imagenes generadas por IA Copilot' paragraph"
}

# --- rewrite that paragraph: split the bold run into "...IA " + "Copilot" ---
# (the second part is wrapped in proofErr spellStart/spellEnd, matching how
# Word flags "Copilot" as an unrecognised word) while keeping the paragraph's
# own numbering / style untouched.
$rewrittenParagraphXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="004807C4" w:rsidRPr="004807C4" w:rsidRDefault="004807C4" w:rsidP="004807C4"><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r w:rsidRPr="004807C4"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="es-MX"/></w:rPr><w:t xml:space="preserve">Imágenes generadas por IA </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="es-MX"/></w:rPr><w:t>Copilot</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>
"@

$targetRange = $target.Range
$targetRange.Collapse(0)
$targetRange.InsertXML($rewrittenParagraphXml)

# --- append the new trailing paragraph ---
$target2 = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Copilot*") {
        $target2 = $para
    }
}
$endRange = $target2.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newParaRange = $newPara.Range

$newParagraphXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="es-MX"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:lang w:val="es-MX"/></w:rPr><w:lastRenderedPageBreak/><w:t>El motivo por el cual escojo este evento, es porque siento que puedo representar mínimamente algo cercano a lo que cada uno de los habitantes de esta ciudad puede llegar a vivir, y no solo puedo conocer de vivencias ajenas a través de libros, historias o películas extranjeras de sucesos que no estuvieron relacionados con nosotros, sino que puedo aprender de las vivencias de mis cercanos y de lo que significó para ellos la llegada del metro a Medellín.</w:t></w:r></w:p>
"@

$newParaRange.InsertXML($newParagraphXml)

Write-Output "Edit applied."
